$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current (before) layout in column A:
#  A1 Text
#  A2 You will be presented...
#  A3 For some images...
#  A4 Some images are special...That is, if the previous...  (single combined sentence)
#  A5 After each response...
#  A6 You can take as long...
#  A7 The experiment will have three blocks...
#  A9 (empty, formatted)
#  A11 (empty, formatted)
#
# Target (after) layout in column A:
#  A1 Text
#  A2 You will be presented...
#  A3 For some images...
#  A4 Some images are special: the correct response is the opposite of the preceding correct response.
#  A5 That is, if the previous correct answer was “Left” then the correct response would now be “Right”.
#  A6 After each response...
#  A7 You can take as long...
#  A8 The experiment will have three blocks...
#  A9  (empty, formatted)
#  A11 (empty, formatted)
#  A13 (empty, formatted)  <- new

# Write rows from the bottom up so we never overwrite a value before it has
# been relocated (column A only ever has one value per row at a time). The
# A4/A5 split is written in A4-then-A5 order so new shared-string entries are
# appended in the same order the diff expects (index 6 then index 7).
$ws.Range("A8").Value = "The experiment will have three blocks, each block will take approximately 10 minutes to complete."
# A8 is a brand new cell (row 8 previously did not exist), so give it the
# same "vertical center" formatting (style index 1) used by the other
# instruction rows.
$ws.Range("A8").VerticalAlignment = -4108
$ws.Range("A7").Value = "You can take as long as you like on each image, but the task will not continue until you press the “Left” or “Right” arrow key."
$ws.Range("A6").Value = "After each response, you will be told whether you got the item correct or incorrect."
$ws.Range("A4").Value = "Some images are special: the correct response is the opposite of the preceding correct response."
$ws.Range("A5").Value = "That is, if the previous correct answer was “Left” then the correct response would now be “Right”."

# Add a new blank (but formatted) row at A13, mirroring the existing blank
# rows (A9, A11), which also just carry the "vertical center" style (s="1")
# with no value.
$ws.Range("A13").VerticalAlignment = -4108

# Update selection to match the post-edit active cell reported in the diff.
$ws.Range("A8").Select() | Out-Null
